$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each of the 3 EC consumption profiles (columns A, B, C) is stored as a
# rolling block of 6 rows: 4 data rows, a subtotal row and a zero row.
# New hourly samples are pushed in at the top of every block, the older
# samples shift down one slot, and the oldest (4th) sample of the block
# is retired.
$groupStarts = @(1, 7, 13, 19, 25, 31, 37, 43, 49, 55)

foreach ($g in $groupStarts) {
    $row1 = $g
    $row2 = $g + 1
    $row3 = $g + 2
    $row4 = $g + 3

    # Capture the two rows that are about to be overwritten (bottom-up).
    $prevRow2 = @($ws.Cells.Item($row2, 1).Value(), $ws.Cells.Item($row2, 2).Value(), $ws.Cells.Item($row2, 3).Value())
    $prevRow3 = @($ws.Cells.Item($row3, 1).Value(), $ws.Cells.Item($row3, 2).Value(), $ws.Cells.Item($row3, 3).Value())

    # Row 4 (oldest kept sample) receives what used to be row 3.
    $ws.Cells.Item($row4, 1).Value = $prevRow3[0]
    $ws.Cells.Item($row4, 2).Value = $prevRow3[1]
    $ws.Cells.Item($row4, 3).Value = $prevRow3[2]

    # Row 3 receives what used to be row 2.
    $ws.Cells.Item($row3, 1).Value = $prevRow2[0]
    $ws.Cells.Item($row3, 2).Value = $prevRow2[1]
    $ws.Cells.Item($row3, 3).Value = $prevRow2[2]
}

# Fresh readings for each profile's two newest rows.
$newRow1 = @{
    1  = @(9235, 7162, 7077)
    7  = @(8467, 6405, 6598)
    13 = @(9386, 6289, 6385)
    19 = @(7933, 4609, 4624)
    25 = @(10566, 7921, 8019)
    31 = @(14563, 9290, 9047)
    37 = @(22846, 13406, 11278)
    43 = @(13791, 13088, 11657)
    49 = @(9795, 9224, 8784)
    55 = @(9961, 8720, 8305)
}
$newRow2 = @{
    1  = @(6668, 2856, 2973)
    7  = @(5513, 2471, 2556)
    13 = @(5370, 2316, 2309)
    19 = @(4364, 1440, 1579)
    25 = @(5873, 2606, 2656)
    31 = @(6580, 3214, 3375)
    37 = @(8977, 3985, 4141)
    43 = @(8247, 4344, 4509)
    49 = @(7327, 3234, 3381)
    55 = @(5831, 2479, 2583)
}

foreach ($g in $groupStarts) {
    $vals1 = $newRow1[$g]
    $ws.Cells.Item($g, 1).Value = $vals1[0]
    $ws.Cells.Item($g, 2).Value = $vals1[1]
    $ws.Cells.Item($g, 3).Value = $vals1[2]

    $vals2 = $newRow2[$g]
    $ws.Cells.Item($g + 1, 1).Value = $vals2[0]
    $ws.Cells.Item($g + 1, 2).Value = $vals2[1]
    $ws.Cells.Item($g + 1, 3).Value = $vals2[2]
}

# PVGenerationFactor / column layout fix: widen columns A:C slightly.
$ws.Columns.Item(1).ColumnWidth = 8.8
$ws.Columns.Item(2).ColumnWidth = 8.8
$ws.Columns.Item(3).ColumnWidth = 8.8
